$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so numeric-looking strings
# like "27.720.61" or "1.0000" are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.720.61'
$ws.Range("E2").Value = '  -1.79%  '

$ws.Range("D3").Value = '1.742.22'
$ws.Range("E3").Value = '  -2.37%  '

$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.53%  '

$ws.Range("D5").Value = '332.33'
$ws.Range("E5").Value = '  -1.22%  '

$ws.Range("D6").Value = '1.0000'
$ws.Range("E6").Value = '  +0.37%  '

$ws.Range("D7").Value = '0.3890'
$ws.Range("E7").Value = '  +1.58%  '

$ws.Range("D8").Value = '0.3360'
$ws.Range("E8").Value = '  -2.32%  '

$ws.Range("D9").Value = '45.40'
$ws.Range("E9").Value = '  -4.58%  '

$ws.Range("D10").Value = '1.095'
$ws.Range("E10").Value = '  -5.45%  '

$ws.Range("D11").Value = '0.07112'
$ws.Range("E11").Value = '  -4.30%  '

$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  +0.55%  '

$ws.Range("D13").Value = '21.73'
$ws.Range("E13").Value = '  -5.84%  '

$ws.Range("D14").Value = '6.052'
$ws.Range("E14").Value = '  -5.71%  '

$ws.Range("D15").Value = '1.741.13'
$ws.Range("E15").Value = '  -2.35%  '

$ws.Range("D16").Value = '6.923'
$ws.Range("E16").Value = '  -2.97%  '

$ws.Range("D17").Value = '0.00001042'
$ws.Range("E17").Value = '  -3.81%  '

$ws.Range("D18").Value = '0.06604'
$ws.Range("E18").Value = '  -0.78%  '

$ws.Range("B19").Value = 'Litecoin'
$ws.Range("C19").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D19").Value = '78.69'
$ws.Range("E19").Value = '  -4.99%  '

$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  +0.54%  '

$ws.Range("D21").Value = '16.62'
$ws.Range("E21").Value = '  -5.01%  '

$ws.Range("D22").Value = '6.148'
$ws.Range("E22").Value = '  -4.37%  '

$ws.Range("D23").Value = '27.736.43'
$ws.Range("E23").Value = '  -1.70%  '

$ws.Range("D24").Value = '11.44'
$ws.Range("E24").Value = '  -5.58%  '

$ws.Range("D25").Value = '2.391'
$ws.Range("E25").Value = '  +0.39%  '

$ws.Range("D26").Value = '153.52'
$ws.Range("E26").Value = '  -0.14%  '

$ws.Range("D27").Value = '19.64'
$ws.Range("E27").Value = '  -6.14%  '

$ws.Range("D28").Value = '2.263'
$ws.Range("E28").Value = '  -6.54%  '

$ws.Range("D29").Value = '1.938.65'
$ws.Range("E29").Value = '  -2.38%  '

$ws.Range("D30").Value = '1.269'
$ws.Range("E30").Value = '  -11.35%  '

$ws.Range("D31").Value = '127.00'
$ws.Range("E31").Value = '  -6.00%  '

$ws.Range("D32").Value = '4.053'
$ws.Range("E32").Value = '  +2.58%  '

$ws.Range("D33").Value = '5.709'
$ws.Range("E33").Value = '  -7.36%  '

$ws.Range("D34").Value = '0.08691'
$ws.Range("E34").Value = '  -1.38%  '

$ws.Range("D35").Value = '11.91'
$ws.Range("E35").Value = '  -6.87%  '

$ws.Range("B36").Value = 'WEMIXTOKEN'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").Value = '1.503'
$ws.Range("E36").Value = '  -0.39%  '

$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").Value = '5.057'
$ws.Range("E37").Value = '  -5.13%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.02241'
$ws.Range("E38").Value = '  -7.72%  '

$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '0.06030'
$ws.Range("E39").Value = '  -4.95%  '

$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = '0.6363'
$ws.Range("E40").Value = '  -7.45%  '

$ws.Range("D41").Value = '0.2071'
$ws.Range("E41").Value = '  -5.07%  '

$ws.Range("E42").Value = '  -4.63%  '

$ws.Range("E43").Value = '  +0.48%  '

$ws.Range("D44").Value = '7.780'
$ws.Range("E44").Value = '  -6.75%  '

$ws.Range("D45").Value = '13.45'
$ws.Range("E45").Value = '  -5.80%  '

$ws.Range("D46").Value = '3.801'
$ws.Range("E46").Value = '  -1.36%  '

$ws.Range("D47").Value = '0.5865'
$ws.Range("E47").Value = '  -7.25%  '

$ws.Range("D48").Value = '124.99'
$ws.Range("E48").Value = '  -5.44%  '

$ws.Range("D49").Value = '1.954'
$ws.Range("E49").Value = '  -6.83%  '

$ws.Range("D50").Value = '0.06912'
$ws.Range("E50").Value = '  -7.34%  '

$ws.Range("D51").Value = '1.137'
$ws.Range("E51").Value = '  -5.74%  '

# Restore the default "Normal" style on column D so no style index is
# left on the cells (matches original workbook which had no explicit
# number format on these data cells).
$ws.Range("D2:D51").Style = "Normal"